$wb = $excel.ActiveWorkbook

# Add the new "ColdStorage" sheet as the last tab (after the current last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ColdStorage"

# Header-ish values carried over in columns A & B (row 1)
$ws.Range("A1").Value = "Boostrix Injection"
$ws.Range("B1").Value = "Boostrix Injection"

# PIN codes for the cold-storage delivery check, stored as text (leading
# apostrophe forces Excel to keep the number-looking value as text, which
# is what produces the quotePrefix cell style)
$ws.Range("C1").Value = "'180001"
$ws.Range("C2").Value = "'600055"
$ws.Range("C3").Value = "'600100"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 23.2222222222222
$ws.Columns.Item(2).ColumnWidth = 18.8888888888889
$ws.Columns.Item(3).ColumnWidth = 14.7777777777778

# Page margins (PageSetup works in points; 1 inch = 72 points)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Touching the header/footer forces the (empty) <headerFooter/> element,
# matching the other sheets in the workbook
$ws.PageSetup.CenterHeader = ""

# Selection left on C7, as in the authored sheet
$ws.Range("C7").Select() | Out-Null
